$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("jhagjas")

$ws.Range("J1").Value = 39.11275243759155
$ws.Range("J2").Value = 38.96059703826904
$ws.Range("J3").Value = 39.64492177963257
$ws.Range("J4").Value = 39.48555731773376
$ws.Range("B5").Value = 2232
$ws.Range("E5").Value = 6
$ws.Range("H5").Value = 99.73106230389959
$ws.Range("I5").Value = 0.003590664272890485
$ws.Range("J5").Value = 40.16272592544556
$ws.Range("B6").Value = 2576
$ws.Range("D6").Value = 2556
$ws.Range("E6").Value = 19
$ws.Range("F6").Value = 14
$ws.Range("G6").Value = 99.45525291828794
$ws.Range("H6").Value = 99.2621359223301
$ws.Range("J6").Value = 41.05736660957336
$ws.Range("J7").Value = 41.75398063659668
$ws.Range("B8").Value = 2139
$ws.Range("E8").Value = 5
$ws.Range("H8").Value = 99.76613657623948
$ws.Range("I8").Value = 0.002810304449648712
$ws.Range("J8").Value = 44.15846085548401
$ws.Range("B9").Value = 1767
$ws.Range("E9").Value = 12
$ws.Range("H9").Value = 99.32049830124575
$ws.Range("I9").Value = 0.01078320090805902
$ws.Range("J9").Value = 45.02007651329041
$ws.Range("J10").Value = 44.39277005195618
$ws.Range("J11").Value = 45.71279048919678
$ws.Range("B12").Value = 2538
$ws.Range("E12").Value = 0
$ws.Range("H12").Value = 100
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 45.5863995552063
$ws.Range("B13").Value = 1795
$ws.Range("E13").Value = 1
$ws.Range("H13").Value = 99.94425863991081
$ws.Range("I13").Value = 0.001114206128133705
$ws.Range("J13").Value = 46.052090883255
$ws.Range("B14").Value = 1878
$ws.Range("E14").Value = 5
$ws.Range("H14").Value = 99.73361747469366
$ws.Range("I14").Value = 0.005854177754124534
$ws.Range("J14").Value = 44.654541015625
$ws.Range("J15").Value = 42.44349765777588
$ws.Range("J16").Value = 44.30206179618835
$ws.Range("J17").Value = 42.85961079597473
$ws.Range("B18").Value = 2279
$ws.Range("E18").Value = 3
$ws.Range("H18").Value = 99.86830553116769
$ws.Range("I18").Value = 0.001318101933216169
$ws.Range("J18").Value = 43.35636258125305
$ws.Range("B19").Value = 1991
$ws.Range("E19").Value = 4
$ws.Range("H19").Value = 99.79899497487438
$ws.Range("I19").Value = 0.002013085052843483
$ws.Range("J19").Value = 42.35549330711365
$ws.Range("J20").Value = 42.43965625762939
$ws.Range("J21").Value = 43.16368007659912
$ws.Range("B22").Value = 1518
$ws.Range("E22").Value = 1
$ws.Range("H22").Value = 99.9340804218853
$ws.Range("I22").Value = 0.0006591957811470006
$ws.Range("J22").Value = 42.60353398323059
$ws.Range("J23").Value = 42.91984677314758
$ws.Range("B24").Value = 2600
$ws.Range("D24").Value = 2597
$ws.Range("F24").Value = 3
$ws.Range("G24").Value = 99.88461538461539
$ws.Range("H24").Value = 99.92304732589457
$ws.Range("I24").Value = 0.001922337562475971
$ws.Range("J24").Value = 43.26388025283813
$ws.Range("B25").Value = 1950
$ws.Range("D25").Value = 1948
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = 14
$ws.Range("G25").Value = 99.28644240570846
$ws.Range("H25").Value = 99.94869163673678
$ws.Range("J25").Value = 42.35196018218994
$ws.Range("J26").Value = 42.78260517120361
$ws.Range("B27").Value = 2765
$ws.Range("D27").Value = 2752
$ws.Range("F27").Value = 226
$ws.Range("G27").Value = 92.41101410342512
$ws.Range("H27").Value = 99.56584659913169
$ws.Range("I27").Value = 0.07989258140315542
$ws.Range("J27").Value = 43.22052311897278
$ws.Range("B28").Value = 2634
$ws.Range("D28").Value = 2633
$ws.Range("F28").Value = 22
$ws.Range("G28").Value = 99.1713747645951
$ws.Range("I28").Value = 0.008283132530120483
$ws.Range("J28").Value = 42.01884341239929
$ws.Range("B29").Value = 2048
$ws.Range("D29").Value = 1851
$ws.Range("E29").Value = 196
$ws.Range("F29").Value = 7
$ws.Range("G29").Value = 99.62325080731969
$ws.Range("H29").Value = 90.42501221299463
$ws.Range("I29").Value = 0.1091984938138784
$ws.Range("J29").Value = 42.05958461761475
$ws.Range("B30").Value = 2943
$ws.Range("E30").Value = 5
$ws.Range("H30").Value = 99.83004758667573
$ws.Range("I30").Value = 0.006772773450728073
$ws.Range("J30").Value = 42.51695394515991
$ws.Range("B31").Value = 2996
$ws.Range("D31").Value = 2995
$ws.Range("F31").Value = 9
$ws.Range("G31").Value = 99.70039946737683
$ws.Range("I31").Value = 0.002995008319467554
$ws.Range("J31").Value = 42.76798915863037
$ws.Range("B32").Value = 2619
$ws.Range("D32").Value = 2617
$ws.Range("E32").Value = 1
$ws.Range("F32").Value = 31
$ws.Range("G32").Value = 98.82930513595166
$ws.Range("H32").Value = 99.96180290297937
$ws.Range("I32").Value = 0.0120800302000755
$ws.Range("J32").Value = 42.06166124343872
$ws.Range("J33").Value = 42.6504852771759
$ws.Range("J34").Value = 43.15852308273315
$ws.Range("B35").Value = 2256
$ws.Range("D35").Value = 2255
$ws.Range("F35").Value = 5
$ws.Range("G35").Value = 99.77876106194691
$ws.Range("I35").Value = 0.00221141088014153
$ws.Range("J35").Value = 40.04944801330566
$ws.Range("B36").Value = 3343
$ws.Range("D36").Value = 3342
$ws.Range("F36").Value = 19
$ws.Range("G36").Value = 99.43469205593573
$ws.Range("I36").Value = 0.005651397977394408
$ws.Range("J36").Value = 40.97624850273132
$ws.Range("B37").Value = 2295
$ws.Range("E37").Value = 91
$ws.Range("H37").Value = 96.03312990409765
$ws.Range("I37").Value = 0.04302536231884058
$ws.Range("J37").Value = 41.48243403434753
$ws.Range("J38").Value = 40.40843796730042
$ws.Range("J39").Value = 40.19192481040955
$ws.Range("J40").Value = 39.76478838920593
$ws.Range("B41").Value = 2480
$ws.Range("D41").Value = 2474
$ws.Range("F41").Value = 7
$ws.Range("G41").Value = 99.71785570334542
$ws.Range("H41").Value = 99.79830576845502
$ws.Range("I41").Value = 0.004834810636583401
$ws.Range("J41").Value = 40.27411818504333
$ws.Range("J42").Value = 39.99424004554749
$ws.Range("B43").Value = 2051
$ws.Range("E43").Value = 6
$ws.Range("H43").Value = 99.70731707317073
$ws.Range("I43").Value = 0.006335282651072124
$ws.Range("J43").Value = 40.10549712181091
$ws.Range("B44").Value = 2255
$ws.Range("E44").Value = 0
$ws.Range("H44").Value = 100
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 40.19294333457947
$ws.Range("B45").Value = 1571
$ws.Range("E45").Value = 1
$ws.Range("H45").Value = 99.93630573248407
$ws.Range("I45").Value = 0.0006369426751592356
$ws.Range("J45").Value = 39.69635796546936
$ws.Range("B46").Value = 1783
$ws.Range("E46").Value = 3
$ws.Range("H46").Value = 99.83164983164983
$ws.Range("I46").Value = 0.001685393258426966
$ws.Range("J46").Value = 39.73055195808411
$ws.Range("B47").Value = 3055
$ws.Range("D47").Value = 3053
$ws.Range("F47").Value = 23
$ws.Range("G47").Value = 99.25227568270481
$ws.Range("H47").Value = 99.96725605762934
$ws.Range("I47").Value = 0.007799805004874878
$ws.Range("J47").Value = 41.0127592086792
$ws.Range("J48").Value = 40.36393451690674

$ws.Name = "sdfewf"

Write-Output "done"
